$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill cells in the exact order the shared strings were authored, so the
# resulting sharedStrings.xml table lines up index-for-index with the target.

# Row 7
$ws.Range("A7").Value = "TC006"
$ws.Range("B7").Value = "Search and Reset Functionality"
$ws.Range("C7").Value = "Verify that the search and reset functionalities work correctly."
$ws.Range("E7").Value = "1. Enter a search term (e.g., First Name or Department).`n2.Verify that only matching records are displayed.`n3. Clear the search input field and verify that the table displays all records.`n4. Test case-insensitivity by searching with uppercase and lowercase terms (e.g., `"JOHN`" vs. `"john`")."
$ws.Range("F7").Value = "The table filters records accurately based on the search term.`nClearing the search resets the table to show all records.`nSearch functionality is case insensitive."
$ws.Range("D7").Value = "Navigate to Web Tables page"
$ws.Range("G7").Value = "Pass"

# Row 8 -- note "Sorting" (B8) was typed before "TC007" (A8) by the author
$ws.Range("B8").Value = "Sorting"
$ws.Range("A8").Value = "TC007"
$ws.Range("C8").Value = "Verify that sorting works correctly."
$ws.Range("E8").Value = "1. Click any column header to apply sorting.`n2. Verify that tab;e is sortedin ascending and descending order for each column."
$ws.Range("F8").Value = "Sorting work as expected and coirrect records are displayed."
$ws.Range("D8").Value = "Navigate to Web Tables page"
$ws.Range("G8").Value = "Pass"

# Row 9
$ws.Range("A9").Value = "TC008"

# Row 6: add "Pass" in column G (G6) -- reuses existing shared string, order
# relative to the new unique strings above does not matter.
$ws.Range("G6").Value = "Pass"

# Apply styles: E7, F7, E8 wrap text (style index 1), G6, G7, G8 center (style index 3)
$ws.Range("E7:F7").WrapText = $true
$ws.Range("E8").WrapText = $true
$ws.Range("G6").HorizontalAlignment = -4108
$ws.Range("G6").VerticalAlignment = -4108
$ws.Range("G7").HorizontalAlignment = -4108
$ws.Range("G7").VerticalAlignment = -4108
$ws.Range("G8").HorizontalAlignment = -4108
$ws.Range("G8").VerticalAlignment = -4108

# row heights (auto-fit result of the wrapped text in the new rows)
$ws.Rows.Item(7).RowHeight = 165
$ws.Rows.Item(8).RowHeight = 75

# view state: zoom + final selection, like the author left the sheet
$excel.ActiveWindow.Zoom = 101
$ws.Range("G9").Select()
